$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header in column L was "availableDateNote"; rename it to "description".
$ws.Range("L1").Value = "description"

# Move/restore the active selection as recorded in the saved file.
$ws.Range("L13").Select()
